# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets,
# plus a single cell on 演出, to match the refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14884
$ws1.Range("F3").Value = 18546
$ws1.Range("F14").Value = 108
$ws1.Range("F15").Value = 200
$ws1.Range("F16").Value = 54
$ws1.Range("F17").Value = 1418
$ws1.Range("F21").Value = 228
$ws1.Range("F22").Value = 7694
$ws1.Range("F23").Value = 985
$ws1.Range("F24").Value = 21
$ws1.Range("F25").Value = 52
$ws1.Range("F26").Value = 1220
$ws1.Range("F28").Value = 5957
$ws1.Range("F29").Value = 101
$ws1.Range("F30").Value = 63
$ws1.Range("F32").Value = 153
$ws1.Range("F33").Value = 257
$ws1.Range("F34").Value = 5310

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14884
$ws4.Range("F3").Value = 18546
$ws4.Range("F14").Value = 108
$ws4.Range("F15").Value = 200
$ws4.Range("F16").Value = 54
$ws4.Range("F17").Value = 1418
$ws4.Range("F22").Value = 228
$ws4.Range("F23").Value = 7694
$ws4.Range("F24").Value = 985
$ws4.Range("F25").Value = 21
$ws4.Range("F26").Value = 52
$ws4.Range("F27").Value = 1220
$ws4.Range("F29").Value = 1
$ws4.Range("F31").Value = 5957
$ws4.Range("F32").Value = 101
$ws4.Range("F33").Value = 63
$ws4.Range("F35").Value = 153
$ws4.Range("F36").Value = 257
$ws4.Range("F37").Value = 5310
